# Insert a new price-report row at row 98 (pushing the existing rows 98-128
# down to 99-129) on the single worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a blank row at position 98.
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new weekly price record.
$ws.Range("A98").Value = 1
$ws.Range("B98").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C98").Value = "Arica y Parinacota"
$ws.Range("D98").Value = 45173
$ws.Range("E98").Value = 15
$ws.Range("F98").Value = 100112040
$ws.Range("G98").Value = "Cilantro"
$ws.Range("H98").Value = "Sin especificar"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 300
$ws.Range("K98").Value = 900
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 950
$ws.Range("N98").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O98").Value = "Región de Arica y Parinacota"
$ws.Range("P98").Value = 475
$ws.Range("Q98").Value = 2
$ws.Range("R98").Value = "Hortaliza"
